$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "This is Version 3 of File"

$ws.Range("C7").Select()
